$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.4131636533673999
$ws.Range("B3").Value = 0.8238509115745072
$ws.Range("B4").Value = 1.400901628879421
$ws.Range("B5").Value = 2.662840662933722
$ws.Range("B6").Value = 3.565271434419023
$ws.Range("B7").Value = 4.390099044047414
$ws.Range("B8").Value = 5.084705158429044
$ws.Range("B9").Value = 6.058480790704488
$ws.Range("B10").Value = 6.480543120354559
$ws.Range("B11").Value = 7.744274164899728
$ws.Range("B12").Value = 8.213118453516353
$ws.Range("B13").Value = 9.771414924822103
$ws.Range("B14").Value = 11.13288767852672
$ws.Range("B15").Value = 11.62438237848504
$ws.Range("B16").Value = 13.23286338686713
$ws.Range("B17").Value = 14.2666822772573
$ws.Range("B18").Value = 15.35775383255386
$ws.Range("B19").Value = 16.55257501767783
$ws.Range("B20").Value = 17.99371871788277
$ws.Range("B21").Value = 19.31498754323765
$ws.Range("B22").Value = 20.10285869029084
$ws.Range("B23").Value = 21.69747103584839
$ws.Range("B24").Value = 23.38256577753764
$ws.Range("B25").Value = 24.87145510509535
$ws.Range("B26").Value = 26.6637355542847
$ws.Range("B27").Value = 27.29120994675154
$ws.Range("B28").Value = 29.46564032812875
$ws.Range("B29").Value = 31.54783256081483
$ws.Range("B30").Value = 32.42340550727129
$ws.Range("B31").Value = 34.72985608586661
$ws.Range("B32").Value = 37.60214223394708
$ws.Range("B33").Value = 39.71323306898522
$ws.Range("B34").Value = 40.9888913641441
$ws.Range("B35").Value = 43.11569601138837
$ws.Range("B36").Value = 45.47121598881147
$ws.Range("B37").Value = 47.24117131621497
$ws.Range("B38").Value = 50.1384477071067
$ws.Range("B39").Value = 52.36126051287686
$ws.Range("B40").Value = 55.03906701522532
$ws.Range("B41").Value = 58.38839046960703
$ws.Range("B42").Value = 60.43854851432661
$ws.Range("B43").Value = 63.56836076088637
$ws.Range("B44").Value = 66.62301075635303
$ws.Range("B45").Value = 69.43257032278152
$ws.Range("B46").Value = 73.26128353070021
$ws.Range("B47").Value = 76.41689312478887
$ws.Range("B48").Value = 79.2507145551281
$ws.Range("B49").Value = 83.25490107051952
$ws.Range("B50").Value = 86.73336770651613
